$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells keep their original text representation (avoid Excel
# auto-converting strings like "1.001" or "30.411.54" into numbers/dates).
$cells = @{
    "D2" = "30.411.54"
    "E2" = "  +1.69%  "
    "D3" = "2.013.12"
    "E3" = "  +4.90%  "
    "E4" = "  +0.03%  "
    "D5" = "325.19"
    "E5" = "  +1.64%  "
    "E6" = "  +0.03%  "
    "D7" = "0.5134"
    "E7" = "  +1.62%  "
    "D8" = "0.4260"
    "E8" = "  +5.22%  "
    "D9" = "0.08754"
    "E9" = "  +4.84%  "
    "D10" = "43.51"
    "E10" = "  +2.63%  "
    "E11" = "  +2.78%  "
    "D12" = "24.44"
    "E12" = "  +2.59%  "
    "D13" = "2.009.54"
    "E13" = "  +5.13%  "
    "D14" = "6.646"
    "E14" = "  +3.69%  "
    "D15" = "7.467"
    "E15" = "  +3.25%  "
    "D16" = "1.001"
    "E16" = "  +0.51%  "
    "D17" = "94.28"
    "E17" = "  +2.30%  "
    "E18" = "  +1.31%  "
    "D19" = "0.06538"
    "E19" = "  +0.32%  "
    "E20" = "  +3.30%  "
    "E21" = "  +0.07%  "
    "D22" = "6.205"
    "E22" = "  +4.33%  "
    "D23" = "30.458.90"
    "E23" = "  +1.72%  "
    "E24" = "  +4.39%  "
    "D25" = "2.254"
    "E25" = "  +2.96%  "
    "D26" = "2.256.54"
    "E26" = "  +5.79%  "
    "E27" = "  +1.85%  "
    "D28" = "162.11"
    "E28" = "  -0.14%  "
    "D29" = "2.447"
    "E29" = "  +4.83%  "
    "D30" = "131.30"
    "E30" = "  +1.91%  "
    "D31" = "1.137"
    "E31" = "  +0.60%  "
    "D32" = "0.1056"
    "E32" = "  +1.82%  "
    "D33" = "6.085"
    "E33" = "  +1.98%  "
    "D34" = "3.831"
    "E34" = "  +1.82%  "
    "D35" = "1.369"
    "E35" = "  +14.67%  "
    "D36" = "0.02534"
    "E36" = "  +3.57%  "
    "D37" = "5.477"
    "E37" = "  +1.20%  "
    "D38" = "0.06672"
    "E38" = "  +4.01%  "
    "B39" = "FraxShare"
    "C39" = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
    "D39" = "9.246"
    "E39" = "  +5.61%  "
    "B40" = "Aptos"
    "C40" = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
    "D40" = "12.40"
    "E40" = "  +9.02%  "
    "D41" = "0.2216"
    "E41" = "  +2.94%  "
    "D42" = "0.6664"
    "E42" = "  +0.81%  "
    "D43" = "1.239"
    "E43" = "  +1.93%  "
    "E44" = "  +0.04%  "
    "D45" = "13.71"
    "E45" = "  +1.34%  "
    "D46" = "0.6180"
    "E46" = "  +1.28%  "
    "D47" = "2.207"
    "E47" = "  -0.96%  "
    "E48" = "  +0.68%  "
    "D49" = "1.260"
    "E49" = "  +4.24%  "
    "D50" = "125.36"
    "E50" = "  +2.94%  "
    "D51" = "81.23"
    "E51" = "  +2.74%  "
}

foreach ($addr in $cells.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $cells[$addr]
}
